$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows with refreshed price (D) / volume-change (E) figures ---
$ws.Range("D2").Value = "70.910.95"
$ws.Range("E2").Value = "  +3.30%  "

$ws.Range("D3").Value = "3.564.32"
$ws.Range("E3").Value = "  +2.29%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.23"
$ws.Range("E5").Value = "  +2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.68"
$ws.Range("E6").Value = "  +2.23%  "

$ws.Range("E7").Value = "  +3.05%  "

$ws.Range("D8").Value = "3.554.86"

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.221"
$ws.Range("E10").Value = "  +21.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.52"
$ws.Range("E12").Value = "  +1.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000318"
$ws.Range("E13").Value = "  +6.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").Value = "4.129.92"
$ws.Range("E15").Value = "  +2.10%  "

$ws.Range("D16").Value = "70.925.91"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.21"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "3.565.87"
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "570.23"
$ws.Range("E20").Value = "  +5.89%  "

$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.64"
$ws.Range("E23").Value = "  -9.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.56"
$ws.Range("E24").Value = "  +4.06%  "

$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.09"
$ws.Range("E26").Value = "  +1.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  +4.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.56"
$ws.Range("E30").Value = "  +4.02%  "

$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.29"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("E33").Value = "  +3.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "551.29"
$ws.Range("E36").Value = "  -3.19%  "

$ws.Range("E37").Value = "  +4.74%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +10.79%  "

$ws.Range("D40").Value = "0.0₃0802"
$ws.Range("E40").Value = "  +5.72%  "

$ws.Range("D42").Value = "3.566.62"
$ws.Range("E42").Value = "  +11.65%  "

$ws.Range("E43").Value = "  +4.85%  "

$ws.Range("E44").Value = "  +3.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0466"
$ws.Range("E45").Value = "  +6.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.33"
$ws.Range("E48").Value = "  +3.45%  "

$ws.Range("E49").Value = "  +3.17%  "

$ws.Range("E50").Value = "  +15.73%  "

$ws.Range("E51").Value = "  +0.05%  "

# --- Rows 34 and 35 swapped (OKB now ranks above Fetch.AI) ---
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.15"
$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.42"
$ws.Range("E35").Value = "  +13.97%  "

